# Updates the Price (D) and Volume(1h) (E) columns of the cryptos sheet
# to the latest scraped values, as committed by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.451.48"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "'3.754.89"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'594.95"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'169.75"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").Value = "'3.752.03"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").Value = "'6.47"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("E13").Value = "  +5.01%  "
$ws.Range("D14").Value = "'36.69"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "'4.388.34"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "'3.764.47"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("D17").Value = "'18.88"
$ws.Range("E17").Value = "  +3.82%  "
$ws.Range("D18").Value = "'67.536.91"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "'7.24"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("E21").Value = "  -4.01%  "
$ws.Range("D22").Value = "'468.95"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").Value = "'0.720"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "'0.0000148"
$ws.Range("E24").Value = "  -6.83%  "
$ws.Range("D25").Value = "'83.81"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'12.17"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "'10.38"
$ws.Range("E28").Value = "  +3.22%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").Value = "'3.908.99"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "'7.68"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("D34").Value = "'30.45"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "'9.13"
$ws.Range("E35").Value = "  -4.65%  "
$ws.Range("D36").Value = "'3.723.71"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("E37").Value = "  +9.01%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").Value = "'0.138"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").Value = "'5.89"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'0.314"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D45").Value = "'8.73"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").Value = "'1.95"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'45.81"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").Value = "'397.77"
$ws.Range("E48").Value = "  -5.37%  "
$ws.Range("D49").Value = "'0.000270"
$ws.Range("E49").Value = "  -6.68%  "
$ws.Range("D50").Value = "'142.03"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -0.37%  "
